$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1983.5454
$ws.Range("I62").Value = 1997.6666
$ws.Range("K62").Value = 1997.6666
$ws.Range("M62").Value = -1373.6666
$ws.Range("H65").Value = 1983.5454
$ws.Range("I65").Value = 1997.6666
$ws.Range("K65").Value = 9988.333000000001
$ws.Range("M65").Value = -6868.333000000001
$ws.Range("H93").Value = 30601
$ws.Range("J93").Value = 30601
$ws.Range("L93").Value = 30601
$ws.Range("N93").Value = -35593
$ws.Range("H103").Value = 825.4167
$ws.Range("J103").Value = 862.2222
$ws.Range("L103").Value = 2586.6666
$ws.Range("N103").Value = -3758.6666
$ws.Range("H129").Value = 2357.0442
$ws.Range("I129").Value = 10477.2
$ws.Range("J129").Value = 957.0172
$ws.Range("K129").Value = 31431.6
$ws.Range("L129").Value = 2871.0516
$ws.Range("M129").Value = -26431.6
$ws.Range("N129").Value = -12871.0516
$ws.Range("H132").Value = 7582379
$ws.Range("I132").Value = 8071472.5
$ws.Range("K132").Value = 24214417.5
$ws.Range("M132").Value = -24211887.5
$ws.Range("H137").Value = 1768.8276
$ws.Range("I137").Value = 1017.8125
$ws.Range("J137").Value = 2693.1538
$ws.Range("K137").Value = 3053.4375
$ws.Range("L137").Value = 8079.4614
$ws.Range("M137").Value = -503.4375
$ws.Range("N137").Value = -13179.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21519.14
$ws.Range("I32").Value = 3589.5693
$ws.Range("J32").Value = 120821.38
$ws.Range("K32").Value = 3589.5693
$ws.Range("L32").Value = 120821.38
$ws.Range("M32").Value = -3302.5693
$ws.Range("N32").Value = -121395.38
$ws.Range("H45").Value = 1673.84
$ws.Range("I45").Value = 1576.3572
$ws.Range("J45").Value = 1797.909
$ws.Range("K45").Value = 1576.3572
$ws.Range("L45").Value = 1797.909
$ws.Range("M45").Value = -1199.3572
$ws.Range("N45").Value = -2551.909
$ws.Range("H61").Value = 1759.3158
$ws.Range("I61").Value = 912.7368
$ws.Range("J61").Value = 2605.8948
$ws.Range("K61").Value = 912.7368
$ws.Range("L61").Value = 2605.8948
$ws.Range("M61").Value = -700.7368
$ws.Range("N61").Value = -3029.8948
$ws.Range("H136").Value = 1759.3158
$ws.Range("I136").Value = 912.7368
$ws.Range("J136").Value = 2605.8948
$ws.Range("K136").Value = 2738.2104
$ws.Range("L136").Value = 7817.6844
$ws.Range("M136").Value = -188.2103999999999
$ws.Range("N136").Value = -12917.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 33000
$ws.Range("J108").Value = 33000
$ws.Range("L108").Value = 33000
$ws.Range("N108").Value = -40680
$ws.Range("H123").Value = 39800
$ws.Range("J123").Value = 39800
$ws.Range("L123").Value = 39800
$ws.Range("N123").Value = -49600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14924.8
$ws.Range("I31").Value = 29434.97
$ws.Range("J31").Value = 2228.4
$ws.Range("K31").Value = 29434.97
$ws.Range("L31").Value = 2228.4
$ws.Range("M31").Value = -29139.97
$ws.Range("N31").Value = -2818.4
$ws.Range("H34").Value = 14924.8
$ws.Range("I34").Value = 29434.97
$ws.Range("J34").Value = 2228.4
$ws.Range("K34").Value = 29434.97
$ws.Range("L34").Value = 2228.4
$ws.Range("M34").Value = -29232.97
$ws.Range("N34").Value = -2632.4
$ws.Range("H45").Value = 8666.666999999999
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 10200
$ws.Range("K45").Value = 1000
$ws.Range("L45").Value = 10200
$ws.Range("M45").Value = -407
$ws.Range("N45").Value = -11386
$ws.Range("H53").Value = 13250
$ws.Range("J53").Value = 13250
$ws.Range("L53").Value = 13250
$ws.Range("N53").Value = -14464
$ws.Range("H58").Value = 11024.37
$ws.Range("I58").Value = 1820.6875
$ws.Range("J58").Value = 24411.545
$ws.Range("K58").Value = 1820.6875
$ws.Range("L58").Value = 24411.545
$ws.Range("M58").Value = -1617.6875
$ws.Range("N58").Value = -24817.545
$ws.Range("H108").Value = 28389
$ws.Range("J108").Value = 28389
$ws.Range("L108").Value = 28389
$ws.Range("N108").Value = -36069
$ws.Range("H122").Value = 2275.2188
$ws.Range("I122").Value = 2199.7307
$ws.Range("K122").Value = 6599.1921
$ws.Range("M122").Value = -4149.1921
$ws.Range("H132").Value = 3015.2632
$ws.Range("I132").Value = 2726.0667
$ws.Range("J132").Value = 4099.75
$ws.Range("K132").Value = 8178.2001
$ws.Range("L132").Value = 12299.25
$ws.Range("M132").Value = -5648.2001
$ws.Range("N132").Value = -17359.25
$ws.Range("H134").Value = 1027.4412
$ws.Range("I134").Value = 952.76666
$ws.Range("J134").Value = 1587.5
$ws.Range("K134").Value = 2858.29998
$ws.Range("L134").Value = 4762.5
$ws.Range("M134").Value = -323.2999799999998
$ws.Range("N134").Value = -9832.5
$ws.Range("H136").Value = 11024.37
$ws.Range("I136").Value = 1820.6875
$ws.Range("J136").Value = 24411.545
$ws.Range("K136").Value = 5462.0625
$ws.Range("L136").Value = 73234.63499999999
$ws.Range("M136").Value = -2912.0625
$ws.Range("N136").Value = -78334.63499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 338.33334
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 386
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 1158
$ws.Range("M29").Value = -23
$ws.Range("N29").Value = -1712
$ws.Range("H46").Value = 250400
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 250400
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 751200
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -751382
$ws.Range("H55").Value = 8160.0557
$ws.Range("J55").Value = 8604.764999999999
$ws.Range("L55").Value = 25814.295
$ws.Range("N55").Value = -26168.295
$ws.Range("H60").Value = 464
$ws.Range("I60").Value = 374.66666
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 1123.99998
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -872.9999800000001
$ws.Range("N60").Value = -3502
$ws.Range("H69").Value = 1664.8182
$ws.Range("I69").Value = 200
$ws.Range("J69").Value = 1811.3
$ws.Range("K69").Value = 600
$ws.Range("L69").Value = 5433.9
$ws.Range("M69").Value = 211
$ws.Range("N69").Value = -7055.9
$ws.Range("H72").Value = 1664.8182
$ws.Range("I72").Value = 200
$ws.Range("J72").Value = 1811.3
$ws.Range("K72").Value = 1800
$ws.Range("L72").Value = 16301.7
$ws.Range("M72").Value = 2256
$ws.Range("N72").Value = -24413.7
$ws.Range("H107").Value = 236298.39
$ws.Range("I107").Value = 360.78946
$ws.Range("K107").Value = 1082.36838
$ws.Range("M107").Value = 837.6316199999999
$ws.Range("H131").Value = 1125.24
$ws.Range("J131").Value = 1145.0947
$ws.Range("L131").Value = 3435.2841
$ws.Range("N131").Value = -13515.2841

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3334.9048
$ws.Range("I122").Value = 3775.9167
$ws.Range("K122").Value = 11327.7501
$ws.Range("M122").Value = -8877.750100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 502490
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 4980
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 4980
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -5252
$ws.Range("H61").Value = 1853.3334
$ws.Range("I61").Value = 1512.5
$ws.Range("K61").Value = 1512.5
$ws.Range("M61").Value = -1310.5
$ws.Range("H100").Value = 2149.75
$ws.Range("I100").Value = 2200
$ws.Range("K100").Value = 2200
$ws.Range("M100").Value = -1659
$ws.Range("H113").Value = 1853.3334
$ws.Range("I113").Value = 1512.5
$ws.Range("K113").Value = 1512.5
$ws.Range("M113").Value = 657.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20452
$ws.Range("H96").Value = 2079.6667
$ws.Range("I96").Value = 2136.625
$ws.Range("K96").Value = 2136.625
$ws.Range("M96").Value = -763.625
$ws.Range("H132").Value = 3251.2703
$ws.Range("I132").Value = 3191.1875
$ws.Range("K132").Value = 9573.5625
$ws.Range("M132").Value = -7043.5625
